$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.493.41'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.596.50'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  -0.52%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '510.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.07'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.17%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -5.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.604.43'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.70'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +9.19%  '
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('E12').Value = '  -1.53%  '
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.052.57'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.482.07'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.69'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.95%  '
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.597.33'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.90%  '
$ws.Range('E19').Value = '  -1.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '350.93'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.56'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.52'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.45'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '152.26'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.75'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.871'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('E38').Value = '  -3.43%  '
$ws.Range('E39').Value = '  -3.63%  '
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.77'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '300.67'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.81%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.101'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.622'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.90'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.77%  '
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.003.70'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.31%  '
